# Akeneo Test Data.xlsx - "akeaneo connector more test"
#
# 1) "Data Read": scroll/select near column I, widen column I.
# 2) "Data to Write Test 1": select F2:F7 (no longer the active tab).
# 3) Insert a new blank worksheet "Sheet2" between "Data to Write Test 1"
#    and "Data to Write ALL".
# 4) "Data to Write ALL": drop the separate "Currency" column, fold the
#    currency into the "Amount" column (now headed "USD", formatted as
#    text, values like "1 USD"), select F8 and make this sheet active.

$wb = $excel.ActiveWorkbook

# --- 1) Data Read -----------------------------------------------------
$wsRead = $wb.Worksheets.Item("Data Read")
$wsRead.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 7
$wsRead.Columns.Item(9).ColumnWidth = 138.25
$wsRead.Range("I8").Select()

# --- 2) Data to Write Test 1 ------------------------------------------
$wsTest1 = $wb.Worksheets.Item("Data to Write Test 1")
$wsTest1.Range("F2:F7").Select()

# --- 3) New blank "Sheet2" between the two "Data to Write" sheets -----
$afterSheet = $wb.Worksheets.Item("Data to Write Test 1")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Sheet2"

# --- 4) Data to Write ALL ----------------------------------------------
$wsAll = $wb.Worksheets.Item("Data to Write ALL")

# Remove the standalone "Currency" column (old column F); everything to
# the right shifts left by one.
$wsAll.Columns.Item(6).Delete()

# The old "Amount" column (now column E) becomes a combined "USD" column
# holding text like "1 USD" instead of a bare number.
$wsAll.Range("E1").Value = "USD"
$wsAll.Range("E1").NumberFormat = "@"

$usdValues = @("1 USD", "2 USD", "3 USD", "4 USD", "5 USD", "6 USD")
for ($i = 0; $i -lt $usdValues.Length; $i++) {
    $cell = $wsAll.Cells.Item($i + 2, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $usdValues[$i]
}

$wsAll.Range("F8").Select()
$wsAll.Activate()
